# Updated cryptos list on Mon May  8 22:22:39 UTC 2023 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns for each coin row, and
# update row 51 from NEARProtocol to RenderToken.
# Numeric-looking price strings are prefixed with a leading apostrophe so
# Excel keeps them as literal text (matching the original inline-string
# cell content) instead of re-parsing/reformatting them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.620.92"
$ws.Range("E2").Value = "  -4.58%  "
$ws.Range("D3").Value = "1.844.83"
$ws.Range("E3").Value = "  -3.97%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").Value = "'312.89"
$ws.Range("E5").Value = "  -3.64%  "
$ws.Range("D6").Value = "'0.9990"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").Value = "'0.4250"
$ws.Range("E7").Value = "  -7.08%  "
$ws.Range("D8").Value = "'0.3631"
$ws.Range("E8").Value = "  -4.56%  "
$ws.Range("D9").Value = "'43.73"
$ws.Range("E9").Value = "  -4.04%  "
$ws.Range("D10").Value = "'0.07219"
$ws.Range("E10").Value = "  -6.78%  "
$ws.Range("D11").Value = "'0.8975"
$ws.Range("E11").Value = "  -8.09%  "
$ws.Range("D12").Value = "'20.69"
$ws.Range("E12").Value = "  -7.54%  "
$ws.Range("D13").Value = "1.831.66"
$ws.Range("E13").Value = "  -4.54%  "
$ws.Range("E14").Value = "  -5.57%  "
$ws.Range("D15").Value = "'5.328"
$ws.Range("E15").Value = "  -6.54%  "
$ws.Range("D16").Value = "'0.06823"
$ws.Range("E16").Value = "  -2.37%  "
$ws.Range("D17").Value = "'1.0000"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "'77.44"
$ws.Range("E18").Value = "  -8.37%  "
$ws.Range("D19").Value = "'0.000008871"
$ws.Range("E19").Value = "  -6.40%  "
$ws.Range("D20").Value = "'0.9991"
$ws.Range("D21").Value = "'15.38"
$ws.Range("E21").Value = "  -7.81%  "
$ws.Range("D22").Value = "27.592.07"
$ws.Range("E22").Value = "  -4.73%  "
$ws.Range("D23").Value = "'4.945"
$ws.Range("E23").Value = "  -7.31%  "
$ws.Range("D24").Value = "'10.74"
$ws.Range("E24").Value = "  -3.11%  "
$ws.Range("D25").Value = "2.066.83"
$ws.Range("E25").Value = "  -4.27%  "
$ws.Range("D26").Value = "'2.046"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "'152.14"
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("D28").Value = "'18.16"
$ws.Range("E28").Value = "  -4.66%  "
$ws.Range("D29").Value = "'5.318"
$ws.Range("E29").Value = "  -5.08%  "
$ws.Range("D30").Value = "'111.08"
$ws.Range("E30").Value = "  -5.50%  "
$ws.Range("D31").Value = "'1.743"
$ws.Range("E31").Value = "  -5.27%  "
$ws.Range("D32").Value = "'0.08886"
$ws.Range("E32").Value = "  -4.45%  "
$ws.Range("D33").Value = "'0.7768"
$ws.Range("E33").Value = "  -10.15%  "
$ws.Range("D34").Value = "'4.478"
$ws.Range("E34").Value = "  -12.25%  "
$ws.Range("D35").Value = "'2.845"
$ws.Range("E35").Value = "  -5.50%  "
$ws.Range("D36").Value = "'1.083"
$ws.Range("E36").Value = "  -12.67%  "
$ws.Range("D37").Value = "'0.9991"
$ws.Range("E37").Value = "  -0.48%  "
$ws.Range("D38").Value = "'0.05424"
$ws.Range("E38").Value = "  -4.60%  "
$ws.Range("D39").Value = "'1.087"
$ws.Range("E39").Value = "  -5.42%  "
$ws.Range("E40").Value = "  -3.69%  "
$ws.Range("D41").Value = "'0.01922"
$ws.Range("E41").Value = "  -6.01%  "
$ws.Range("D42").Value = "'0.5045"
$ws.Range("E42").Value = "  -8.17%  "
$ws.Range("D43").Value = "'6.778"
$ws.Range("E43").Value = "  -9.10%  "
$ws.Range("D44").Value = "'0.1633"
$ws.Range("E44").Value = "  -6.95%  "
$ws.Range("D45").Value = "'0.06610"
$ws.Range("E45").Value = "  -4.69%  "
$ws.Range("D46").Value = "'8.235"
$ws.Range("E46").Value = "  -11.71%  "
$ws.Range("D47").Value = "'106.33"
$ws.Range("E47").Value = "  -3.85%  "
$ws.Range("D48").Value = "'0.4712"
$ws.Range("E48").Value = "  -8.70%  "
$ws.Range("D49").Value = "'10.26"
$ws.Range("E49").Value = "  -8.11%  "
$ws.Range("D50").Value = "'0.9983"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'1.876"
$ws.Range("E51").Value = "  -13.68%  "
